$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing data occupies rows 2-21 (20 rows x 3 cols).
# Final data must occupy rows 2-31 (30 rows x 3 cols):
#   - 2 brand-new rows inserted at the top (new rows 2-3)
#   - the old 20 rows shifted down by 2 (now rows 4-23)
#   - 8 brand-new rows appended at the bottom (new rows 24-31)

# Shift the existing 20 rows of data down by two rows, working from the
# bottom up so we don't overwrite values before they are copied.
for ($r = 21; $r -ge 2; $r--) {
    $srcRow = $ws.Range("A" + $r + ":C" + $r)
    $dstRow = $r + 2
    $ws.Range("A" + $dstRow + ":C" + $dstRow).Value2 = $srcRow.Value2
}

# Fill in the two new rows at the top (rows 2-3).
$ws.Range("A2").Value = 0.2804546356201172
$ws.Range("B2").Value = 0.4303635954856872
$ws.Range("C2").Value = -0.691750168800354

$ws.Range("A3").Value = 0.1987819671630859
$ws.Range("B3").Value = 0.2879692316055298
$ws.Range("C3").Value = -0.9282988905906676

# Append eight new rows at the bottom (rows 24-31).
$newRows = @(
    @(-0.4514303207397461, -0.07753515243530271, -1.056098580360413),
    @(1.037992477416992, -1.273390769958496, 0.4362349510192871),
    @(0.0754270553588867, 1.646718859672546, 1.695090532302856),
    @(-0.2560558319091797, 0.3026316165924072, -0.4233262538909912),
    @(0.6335611343383789, 0.8106564879417419, -1.443797469139099),
    @(0.09285736083984369, 0.7357764840126038, -1.646607518196106),
    @(0.0882749557495117, 0.1726978719234466, -0.9354652166366576),
    @(0.2656211853027344, 0.4902379512786865, -0.8409426212310791)
)

$startRow = 24
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($r, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($r, 3).Value = $newRows[$i][2]
}
